# New weekly price observation is published: insert a fresh row at the top
# of the data table (row 2), pushing every existing record down by one row
# (old row 2 -> row 3, old row 3 -> row 4, ... old row 24 -> row 25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2-24) down by one to make room for the new
# record; Excel's InsertRow semantics push row 2..24 -> 3..25 and expand the
# used range / dimension automatically.
$ws.Rows.Item(2).Insert()

# The inserted row picks up formatting from the header row above it (bold,
# centered). The data rows are unformatted (default style) apart from the
# Fecha column, so strip that back off before writing the new values.
$ws.Range("A2:R2").ClearFormats()

# Populate the newly inserted row 2 with this week's observation.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44860
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112003
$ws.Range("G2").Value = "Ajo"
$ws.Range("H2").Value = "Chino"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14500
$ws.Range("N2").Value = "$/caja 10 kilos"
$ws.Range("O2").Value = "China"
$ws.Range("P2").Value = 1450
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = "Hortaliza"

# The Fecha column (D) carries a custom date number format; make sure the
# freshly written cell matches the formatting used by the rest of the
# column (row 3 now holds what used to be row 2, so it still has it).
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat
